# Append one new row (row 9) to Sheet1, matching the source data:
# أحمد شريم | 222 | الصمود | الرحلة 1 | C3 | NRC | ٠١‏/٠٥‏/٢٠٢٥ ٠٧:٤٠:٠١ م
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("ملاحظات") is blank for this row, same as every other row above it.
# A bare "" assignment clears/removes the cell instead of leaving an empty
# text value behind, so seed it with a lone apostrophe (forces text type,
# empty display text) and then strip the resulting quote-prefix style back
# to Normal so it lines up with the rest of the sheet (which carries no
# explicit per-cell style).
$ws.Range("A9").Value = "'"
$ws.Range("A9").Style = "Normal"

$ws.Range("B9").Value = "أحمد شريم"

# Column C ("الكمية") holds numbers-as-text throughout the sheet (see the
# numberStoredAsText ignoredError covering the table). Force "222" to stay
# text the same way, then reset the style so no stray quote-prefix format
# is left on the cell.
$ws.Range("C9").Value = "'222"
$ws.Range("C9").Style = "Normal"

$ws.Range("D9").Value = "الصمود"
$ws.Range("E9").Value = "الرحلة 1"
$ws.Range("F9").Value = "C3"
$ws.Range("G9").Value = "NRC"
$ws.Range("H9").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٧:٤٠:٠١ م"
